$d = $word.ActiveDocument

# --- Paragraphs where Word re-split runs around proofErr marks (spell/grammar check) ---
# Text content is unchanged; only run boundaries + proofErr markers differ.

$d.Paragraphs(2).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Conditional Statements -&gt; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>if</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>,else</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">, else </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>if,switch</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
$d.Paragraphs(4).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> number = 1;</w:t></w:r></w:p>')
$d.Paragraphs(9).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Post </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>increment(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>number++)</w:t></w:r></w:p>')
$d.Paragraphs(10).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Logical </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Operators(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> &amp;&amp;, ||, ==, !=)</w:t></w:r></w:p>')
$d.Paragraphs(18).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>For(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r><w:t>initialization;condition;increment</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> or decrement) {</w:t></w:r></w:p>')
$d.Paragraphs(21).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>For(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">=1;  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> &lt; 11;  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>++ ) {</w:t></w:r></w:p>')
$d.Paragraphs(22).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>System.out.println</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>);// 1,2,3,4,5,6,7,8,9,10</w:t></w:r></w:p>')
$d.Paragraphs(27).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Examples on </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>for  loop</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')
$d.Paragraphs(32).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>WAP  to</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> find given number  is palindrome or not</w:t></w:r></w:p>')
$d.Paragraphs(33).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">121 -&gt; </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">121  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>polindrome</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>')
$d.Paragraphs(36).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>12 -</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>&gt;  1</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>+2+3+4+6</w:t></w:r></w:p>')
$d.Paragraphs(38).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>15 -</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>&gt;  1</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>+3+5 =</w:t></w:r></w:p>')
$d.Paragraphs(44).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">153 -&gt; 1*1*1+5*5*5+3*3*3 </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>=  153</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

# --- Append new content block after paragraph 45 ("125 -> 1*1*1+2*2*2+5*5*5 = 134") ---
$anchor = $d.Paragraphs(45)
$anchor.Range.InsertParagraphAfter()
$newBlock = $d.Paragraphs(46)
$newBlock.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pBdr><w:bottom w:val="double" w:sz="6" w:space="1" w:color="auto"/></w:pBdr></w:pPr><w:r><w:t>24-12-2021</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Triangle Stars Example</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pBdr><w:bottom w:val="double" w:sz="6" w:space="1" w:color="auto"/></w:pBdr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>While Loop</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Syntax:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">       </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Step  1</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> -&gt; true it will go to loop otherwise it will exit from the loop</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>While(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>condition) {</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Step 2</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>}</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 1;</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>While(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> &lt;=10)  {</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>System.out.println</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>);</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>++;</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>}</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

Write-Output "done"